$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.827.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.226.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: Solana -> Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.67%  "

# Row 6: BNB -> BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "629.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.66%  "

# Row 7: Dogecoin -> Dogecoin
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.67%  "

# Row 8: XRP -> XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.698"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.21%  "

# Row 9: USDC -> USDC
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.08%  "

# Row 10: LidoStakedEther -> LidoStakedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.221.19"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.06%  "

# Row 11: Cardano -> Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.572"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.86%  "

# Row 12: TRON -> TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.179"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.08%  "

# Row 13: ShibaInu -> ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.65%  "

# Row 14: Toncoin -> Toncoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.39"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.05%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.822.46"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.02%  "

# Row 16: Avalanche -> Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17: WrappedBTC -> WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.670.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.54%  "

# Row 18: WrappedEther -> WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.212.94"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19: PEPE -> PEPE
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000231"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +72.78%  "

# Row 20: SuiNetwork -> SuiNetwork
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +13.00%  "

# Row 21: BitcoinCash -> Chainlink
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22: Chainlink -> BitcoinCash
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "435.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.55%  "

# Row 23: Uniswap -> Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.64"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.28%  "

# Row 24: Polkadot -> Polkadot
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.28%  "

# Row 25: NEARProtocol -> NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.48%  "

# Row 26: Aptos -> Aptos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.47%  "

# Row 27: Litecoin -> Litecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "79.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.94%  "

# Row 28: WrappedeETH -> WrappedeETH
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.381.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.29%  "

# Row 29: Dai -> Dai
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "

# Row 30: Binance-PegBSC-USD -> Binance-PegBSC-USD
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31: Cronos -> Cronos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.158"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -8.51%  "

# Row 32: dogwifhat -> dogwifhat
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +34.15%  "

# Row 33: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.51"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.45%  "

# Row 34: Bittensor -> Bittensor
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "541.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.67%  "

# Row 35: PancakeSwap -> RenderToken
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.48%  "

# Row 36: RenderToken -> PancakeSwap
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.91"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.66%  "

# Row 37: Fetch.AI -> Fetch.AI
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38: EthereumClassic -> EthereumClassic
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.45"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.46%  "

# Row 39: WhiteBITCoin -> WhiteBITCoin
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.34"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.18%  "

# Row 40: FirstDigitalUSD -> Kaspa
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.127"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.25%  "

# Row 41: Kaspa -> FirstDigitalUSD
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42: USDe -> USDe
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.06%  "

# Row 43: Stacks -> Stacks
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.18%  "

# Row 44: PolygonEcosystemToken -> PolygonEcosystemToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.374"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.02%  "

# Row 45: Monero -> Monero
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.32%  "

# Row 46: Aave -> Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "173.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.74%  "

# Row 47: OKB -> OKB
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.11%  "

# Row 48: Stellar -> Mantle
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.745"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.64%  "

# Row 49: Mantle -> ImmutableX
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.38%  "

# Row 50: ImmutableX -> Stellar
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.123"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.25%  "

# Row 51: ARBITRUM -> ARBITRUM
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.620"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.18%  "

